$d = $word.ActiveDocument

# XML namespace header reused for all InsertXML fragments (Word "single file package" form).
$pkgHeader = "<?xml version='1.0' encoding='UTF-8' standalone='yes'?><pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'><pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'><pkg:xmlData><w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:body>"
$pkgFooter = "</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"

function Remove-WholeParagraphByText([string]$anchorText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Output "NOT FOUND: $anchorText"
        return
    }
    $para = $rng.Paragraphs(1)
    $pr = $para.Range
    $pr.Delete()
}

# 1) Drop the whole "Transformer für IAAS-Referenzen" paragraph - its neighbour
#    ("Globale Einstellungen speichern und abrufen") simply slides up to take its place.
Remove-WholeParagraphByText("Transformer für IAAS-Referenzen")

# 2) Drop the whole "Verwaltung der internen Datenbank" paragraph.
Remove-WholeParagraphByText("Verwaltung der internen Datenbank")

# 3) Remove the stray <w:lastRenderedPageBreak/> that used to sit in front of
#    "Admin-Konsole speichern usw." (the component list was split into its own file,
#    so this break no longer belongs here).
$rngAdmin = $d.Content
$found = $rngAdmin.Find.Execute("Admin-Konsole speichern usw.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $paraAdmin = $rngAdmin.Paragraphs(1)
    $prAdmin = $paraAdmin.Range
    $rNoMark = $d.Range($prAdmin.Start, $prAdmin.End - 1)
    $fragAdmin = $pkgHeader + "<w:p><w:r><w:rPr><w:rFonts w:ascii='Arial' w:hAnsi='Arial' w:cs='Arial'/></w:rPr><w:t>Admin-Konsole speichern usw.</w:t></w:r></w:p>" + $pkgFooter
    $rNoMark.InsertXML($fragAdmin)
} else {
    Write-Output "NOT FOUND: Admin-Konsole speichern usw."
}

# 4) Drop the whole "Referenzen (IAAS) für RRS speichern und abrufen" paragraph.
Remove-WholeParagraphByText("Referenzen (IAAS) für RRS speichern und abrufen")

# 5) The page break now renders one paragraph later, right before the
#    "Externe Funktionen" heading - add it there as its own leading run.
$rngExterne = $d.Content
$found = $rngExterne.Find.Execute("Externe Funktionen", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $paraExterne = $rngExterne.Paragraphs(1)
    $prExterne = $paraExterne.Range
    $rInsert = $d.Range($prExterne.Start, $prExterne.Start)
    $fragBreak = $pkgHeader + "<w:p><w:r><w:rPr><w:rFonts w:ascii='Arial' w:hAnsi='Arial' w:cs='Arial'/></w:rPr><w:lastRenderedPageBreak/></w:r></w:p>" + $pkgFooter
    $rInsert.InsertXML($fragBreak)
} else {
    Write-Output "NOT FOUND: Externe Funktionen"
}

# 6) Remove the <w:lastRenderedPageBreak/> in front of "SIMPL Datasource Service"...
$rngDS = $d.Content
$found = $rngDS.Find.Execute("SIMPL Datasource Service", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $paraDS = $rngDS.Paragraphs(1)
    $prDS = $paraDS.Range
    $rFirst = $d.Range($prDS.Start, $prDS.Start + 1)
    $fragTab = $pkgHeader + "<w:p><w:r w:rsidRPr='004524EA'><w:rPr><w:rFonts w:ascii='Arial' w:hAnsi='Arial' w:cs='Arial'/></w:rPr><w:tab/></w:r></w:p>" + $pkgFooter
    $rFirst.InsertXML($fragTab)
} else {
    Write-Output "NOT FOUND: SIMPL Datasource Service"
}

# 7) ...because it now belongs one paragraph further down, right before
#    "SIMPL Registry Service".
$rngReg = $d.Content
$found = $rngReg.Find.Execute("SIMPL Registry Service", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $paraReg = $rngReg.Paragraphs(1)
    $prReg = $paraReg.Range
    $rNoMarkReg = $d.Range($prReg.Start, $prReg.End - 1)
    $fragReg = $pkgHeader + "<w:p><w:r w:rsidRPr='00E73CFE'><w:rPr><w:rFonts w:ascii='Arial' w:hAnsi='Arial' w:cs='Arial'/><w:lang w:val='en-US'/></w:rPr><w:lastRenderedPageBreak/><w:t>SIMPL Registry Service</w:t></w:r></w:p>" + $pkgFooter
    $rNoMarkReg.InsertXML($fragReg)
} else {
    Write-Output "NOT FOUND: SIMPL Registry Service"
}

Write-Output "done"
